$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some contain multiple "." separators, e.g.
# "62.916.30", which Excel always treats as text). Others look like plain
# decimal numbers (e.g. "528.12") and Excel would silently convert those to
# real numeric values on assignment. To keep every updated D-cell as plain
# text (matching the original inlineStr cells), force the Text number format
# right before writing the value, then restore the default "Normal" style so
# no stray formatting remains on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.155.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.99%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.229.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.28%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.91%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.595'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.52%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.224.11'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.606'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.30'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.134'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.740.62'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.116'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.220.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.913.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.40%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.968'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '366.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.46%  '

$ws.Range("E28").Value = '  +1.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.52'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '636.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.92%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.106'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '56.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.88%  '

$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.376'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0711'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.32%  '

$ws.Range("E42").Value = '  +1.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.878.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0393'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.02'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.47%  '

$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '134.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.79%  '

